$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Type" column between "Code" (B) and "Date" (C), pushing "Date" to D.
# Inserting the column causes it to inherit the adjacent column's style.
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "Type"

# Move the active selection to C1 (as per the diff's selection change)
$ws.Range("C1").Select()
